# Correção das notas do fórum para matc65 em 2021.2
# Zera os valores de presença diária (colunas B:H), total_views (I) e nota_view (J)
# para todas as linhas de dados (2 a 50), mantendo a coluna A (matricula) e o
# cabeçalho (linha 1) intactos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:J50").Value = 0
